# Swap the data values between row 3 and row 4 for the columns that differ
# (A, B, E, F, G, H, Q, R, AC), leaving identical columns (D, P, Q-shared, S,
# T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"

    $val3 = $ws.Range($addr3).Value()
    $val4 = $ws.Range($addr4).Value()

    $ws.Range($addr3).Value = $val4
    $ws.Range($addr4).Value = $val3
}
